# The glossary export re-orders the existing term rows (rows 2-17 of
# Sheet1, columns A-G) into a new sequence; the set of terms/definitions
# itself is unchanged. Implement this as: snapshot every existing row's
# values, then write those values back out in the new row order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Maps each NEW row number -> the OLD row number whose content should end
# up there (derived from the canonical OOXML diff).
$mapping = @{
  2  = 3
  3  = 10
  4  = 2
  5  = 9
  6  = 7
  7  = 8
  8  = 13
  9  = 6
  10 = 12
  11 = 14
  12 = 5
  13 = 15
  14 = 11
  15 = 4
  16 = 17
  17 = 16
}

# 1) Snapshot all current rows (term..sources, columns A-G) before any
#    writes happen, so overlapping source/target rows don't clobber data
#    that still needs to be read.
$oldData = @{}
for ($r = 2; $r -le 17; $r++) {
  $row = @()
  for ($c = 1; $c -le 7; $c++) {
    $row += ,($ws.Cells.Item($r, $c).Value())
  }
  $oldData[$r] = $row
}

# 2) Write the snapshotted rows back out into their new positions.
foreach ($newR in $mapping.Keys) {
  $oldR = $mapping[$newR]
  $row = $oldData[$oldR]
  for ($c = 1; $c -le 7; $c++) {
    $ws.Cells.Item($newR, $c).Value = $row[$c - 1]
  }
}
